$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 54570.45
$ws.Range("I62").Value = 64369.312
$ws.Range("J62").Value = 15375
$ws.Range("K62").Value = 64369.312
$ws.Range("L62").Value = 15375
$ws.Range("M62").Value = -63745.312
$ws.Range("N62").Value = -16623
$ws.Range("H65").Value = 54570.45
$ws.Range("I65").Value = 64369.312
$ws.Range("J65").Value = 15375
$ws.Range("K65").Value = 321846.56
$ws.Range("L65").Value = 76875
$ws.Range("M65").Value = -318726.56
$ws.Range("N65").Value = -83115
$ws.Range("H75").Value = 38804
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 38804
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 38804
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -40676
$ws.Range("H78").Value = 38804
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 38804
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 116412
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -125772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1598.8667
$ws.Range("I2").Value = 1222.5
$ws.Range("J2").Value = 2029
$ws.Range("K2").Value = 1222.5
$ws.Range("L2").Value = 2029
$ws.Range("M2").Value = -1109.5
$ws.Range("N2").Value = -2255
$ws.Range("H61").Value = 9617427
$ws.Range("I61").Value = 15627078
$ws.Range("J61").Value = 1985.6
$ws.Range("K61").Value = 15627078
$ws.Range("L61").Value = 1985.6
$ws.Range("M61").Value = -15626866
$ws.Range("N61").Value = -2409.6
$ws.Range("H88").Value = 3176.2632
$ws.Range("I88").Value = 3109.0908
$ws.Range("J88").Value = 3268.625
$ws.Range("K88").Value = 3109.0908
$ws.Range("L88").Value = 3268.625
$ws.Range("M88").Value = -2703.0908
$ws.Range("N88").Value = -4080.625
$ws.Range("H91").Value = 3176.2632
$ws.Range("I91").Value = 3109.0908
$ws.Range("J91").Value = 3268.625
$ws.Range("K91").Value = 3109.0908
$ws.Range("L91").Value = 3268.625
$ws.Range("M91").Value = -1705.0908
$ws.Range("N91").Value = -6076.625
$ws.Range("H116").Value = 1598.8667
$ws.Range("I116").Value = 1222.5
$ws.Range("J116").Value = 2029
$ws.Range("K116").Value = 1222.5
$ws.Range("L116").Value = 2029
$ws.Range("M116").Value = 1071.5
$ws.Range("N116").Value = -6617
$ws.Range("H122").Value = 2316.6667
$ws.Range("I122").Value = 1980
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5940
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3490
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 4526711
$ws.Range("I132").Value = 1291.2858
$ws.Range("J132").Value = 9806368
$ws.Range("K132").Value = 3873.8574
$ws.Range("L132").Value = 29419104
$ws.Range("M132").Value = -1343.8574
$ws.Range("N132").Value = -29424164
$ws.Range("H133").Value = 31600
$ws.Range("I133").Value = 10000
$ws.Range("J133").Value = 37000
$ws.Range("K133").Value = 10000
$ws.Range("L133").Value = 37000
$ws.Range("M133").Value = -7470
$ws.Range("N133").Value = -42060
$ws.Range("H136").Value = 9617427
$ws.Range("I136").Value = 15627078
$ws.Range("J136").Value = 1985.6
$ws.Range("K136").Value = 46881234
$ws.Range("L136").Value = 5956.799999999999
$ws.Range("M136").Value = -46878684
$ws.Range("N136").Value = -11056.8
$ws.Range("H137").Value = 5980
$ws.Range("I137").Value = 5980
$ws.Range("K137").Value = 5980
$ws.Range("M137").Value = -880
$ws.Range("H139").Value = 58357
$ws.Range("J139").Value = 58357
$ws.Range("L139").Value = 58357
$ws.Range("N139").Value = -68637

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1598.8667
$ws.Range("I3").Value = 1222.5
$ws.Range("J3").Value = 2029
$ws.Range("K3").Value = 1222.5
$ws.Range("L3").Value = 2029
$ws.Range("M3").Value = -1108.5
$ws.Range("N3").Value = -2257
$ws.Range("H86").Value = 1294955.1
$ws.Range("I86").Value = 3392.4546
$ws.Range("J86").Value = 3324553.8
$ws.Range("K86").Value = 3392.4546
$ws.Range("L86").Value = 3324553.8
$ws.Range("M86").Value = -2269.4546
$ws.Range("N86").Value = -3326799.8
$ws.Range("H89").Value = 1294955.1
$ws.Range("I89").Value = 3392.4546
$ws.Range("J89").Value = 3324553.8
$ws.Range("K89").Value = 16962.273
$ws.Range("L89").Value = 16622769
$ws.Range("M89").Value = -11346.273
$ws.Range("N89").Value = -16634001
$ws.Range("H107").Value = 12196572
$ws.Range("I107").Value = 17242358
$ws.Range("K107").Value = 17242358
$ws.Range("M107").Value = -17240438

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4596.5713
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 5355.2
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 5355.2
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -6603.2
$ws.Range("H65").Value = 4596.5713
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 5355.2
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 26776
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -33016
$ws.Range("H134").Value = 15625990
$ws.Range("I134").Value = 995.9231
$ws.Range("J134").Value = 83334300
$ws.Range("K134").Value = 2987.7693
$ws.Range("L134").Value = 250002900
$ws.Range("M134").Value = -452.7692999999999
$ws.Range("N134").Value = -250007970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H116").Value = 166667520
$ws.Range("I116").Value = 866.6667
$ws.Range("J116").Value = 333334180
$ws.Range("K116").Value = 2600.0001
$ws.Range("L116").Value = 1000002540
$ws.Range("M116").Value = 841.9998999999998
$ws.Range("N116").Value = -1000009424
$ws.Range("H131").Value = 864.4536000000001
$ws.Range("J131").Value = 873.6383
$ws.Range("L131").Value = 2620.9149
$ws.Range("N131").Value = -12700.9149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1975
$ws.Range("I113").Value = 950
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 950
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1220
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 42857868
$ws.Range("I16").Value = 3969021.5
$ws.Range("J16").Value = 142857760
$ws.Range("K16").Value = 3969021.5
$ws.Range("L16").Value = 142857760
$ws.Range("M16").Value = -3968851.5
$ws.Range("N16").Value = -142858100
$ws.Range("H61").Value = 4653.5
$ws.Range("I61").Value = 1104
$ws.Range("J61").Value = 5836.6665
$ws.Range("K61").Value = 1104
$ws.Range("L61").Value = 5836.6665
$ws.Range("M61").Value = -902
$ws.Range("N61").Value = -6240.6665
$ws.Range("H93").Value = 1336.9412
$ws.Range("I93").Value = 1226
$ws.Range("K93").Value = 1226
$ws.Range("M93").Value = 22
$ws.Range("H113").Value = 4653.5
$ws.Range("I113").Value = 1104
$ws.Range("J113").Value = 5836.6665
$ws.Range("K113").Value = 1104
$ws.Range("L113").Value = 5836.6665
$ws.Range("M113").Value = 1066
$ws.Range("N113").Value = -10176.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 29847.666
$ws.Range("I122").Value = 51240.8
$ws.Range("J122").Value = 3106.25
$ws.Range("K122").Value = 153722.4
$ws.Range("L122").Value = 9318.75
$ws.Range("M122").Value = -151272.4
$ws.Range("N122").Value = -14218.75
$ws.Range("H130").Value = 59400
$ws.Range("J130").Value = 59400
$ws.Range("L130").Value = 59400
$ws.Range("N130").Value = -69440
